$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells keep their original "General Text" cell type (the source workbook
# stores these as inline/shared strings, not numbers), so force text format on the
# Price (D) column before assigning values that look numeric (e.g. "562.09", "0.348")
# to avoid Excel auto-converting them into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.869.55"
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.414.79"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.09"
$ws.Range("E5").Value = "  +0.99%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.28"
$ws.Range("E6").Value = "  -0.22%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("E9").Value = "  +0.17%  "

$ws.Range("E11").Value = "  -3.75%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "25.68"
$ws.Range("E13").Value = "  -1.42%  "

$ws.Range("E14").Value = "  -1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.848.76"
$ws.Range("E15").Value = "  +0.27%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.780.17"
$ws.Range("E16").Value = "  -0.29%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.424.70"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("E18").Value = "  +1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "323.00"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.78"
$ws.Range("E20").Value = "  +0.79%  "

$ws.Range("E21").Value = "  -1.77%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.44"
$ws.Range("E23").Value = "  +1.74%  "

$ws.Range("E24").Value = "  -0.17%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.74"
$ws.Range("E25").Value = "  -4.70%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "553.28"
$ws.Range("E26").Value = "  -4.66%  "

$ws.Range("E27").Value = "  +0.24%  "

$ws.Range("E28").Value = "  +0.00%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0928"
$ws.Range("E29").Value = "  -0.23%  "

$ws.Range("E30").Value = "  -1.45%  "

$ws.Range("E31").Value = "  -4.67%  "

$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("E33").Value = "  -0.86%  "

$ws.Range("E34").Value = "  -4.50%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.09%  "

$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("E37").Value = "  -1.30%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "153.78"
$ws.Range("E38").Value = "  +2.04%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.39"

$ws.Range("E40").Value = "  -1.02%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.80"
$ws.Range("E41").Value = "  -1.88%  "

$ws.Range("E42").Value = "  -0.87%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "146.55"
$ws.Range("E43").Value = "  -2.94%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.22"
$ws.Range("E44").Value = "  -3.71%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.61"
$ws.Range("E45").Value = "  -0.86%  "

$ws.Range("E46").Value = "  -2.97%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.590"
$ws.Range("E47").Value = "  +0.39%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.69"
$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0917"
$ws.Range("E49").Value = "  -0.56%  "

$ws.Range("E50").Value = "  -0.93%  "
